$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '37.674.16'
$ws.Range("E2").Value = '  -0.15%  '
$ws.Range("D3").Value = '2.037.43'
$ws.Range("E3").Value = '  +0.76%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '227.16'
$ws.Range("E5").Value = '  +0.36%  '
$ws.Range("E6").Value = '  -0.60%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '59.47'
$ws.Range("E7").Value = '  +0.16%  '
$ws.Range("E8").Value = '  -0.05%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.376'
$ws.Range("E9").Value = '  -1.65%  '
$ws.Range("E10").Value = '  +3.49%  '
$ws.Range("E11").Value = '  -0.06%  '
$ws.Range("D12").Value = '2.338.82'
$ws.Range("E12").Value = '  +0.70%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '14.41'
$ws.Range("E13").Value = '  -0.38%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '21.07'
$ws.Range("E14").Value = '  +0.77%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.46'
$ws.Range("E15").Value = '  +4.98%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.769'
$ws.Range("E16").Value = '  +2.96%  '
$ws.Range("D17").Value = '2.033.32'
$ws.Range("E17").Value = '  -0.07%  '
$ws.Range("D18").Value = '37.634.55'
$ws.Range("E18").Value = '  -0.17%  '
$ws.Range("E19").Value = '  -0.98%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '69.27'
$ws.Range("E20").Value = '  -0.02%  '
$ws.Range("E21").Value = '  +0.45%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '223.77'
$ws.Range("E22").Value = '  -0.09%  '
$ws.Range("E23").Value = '  +0.08%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.39'
$ws.Range("E24").Value = '  -0.85%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.27'
$ws.Range("E25").Value = '  +3.30%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.38'
$ws.Range("E26").Value = '  +2.81%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '167.86'
$ws.Range("E27").Value = '  +1.68%  '
$ws.Range("E28").Value = '  -0.60%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '18.77'
$ws.Range("E29").Value = '  -0.09%  '
$ws.Range("E30").Value = '  +0.10%  '
$ws.Range("E31").Value = '  +0.76%  '
$ws.Range("E32").Value = '  +9.75%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.38'
$ws.Range("E33").Value = '  -0.91%  '
$ws.Range("E34").Value = '  +1.56%  '
$ws.Range("E35").Value = '  -0.09%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '6.53'
$ws.Range("E36").Value = '  +3.74%  '
$ws.Range("E37").Value = '  +4.84%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.41'
$ws.Range("E38").Value = '  +5.64%  '
$ws.Range("E39").Value = '  -0.10%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '18.10'
$ws.Range("E40").Value = '  +9.57%  '
$ws.Range("D41").Value = '1.532.82'
$ws.Range("E41").Value = '  +0.05%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '96.95'
$ws.Range("E42").Value = '  +1.33%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0215'
$ws.Range("E43").Value = '  -0.12%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '4.29'
$ws.Range("E45").Value = '  +10.80%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0909'
$ws.Range("E46").Value = '  -0.76%  '
$ws.Range("E47").Value = '  +0.68%  '
$ws.Range("E48").Value = '  +1.00%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.94'
$ws.Range("E49").Value = '  -0.52%  '
$ws.Range("E50").Value = '  -0.10%  '
$ws.Range("D51").Value = '2.228.03'
$ws.Range("E51").Value = '  +0.69%  '
